$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D1").Value = "Test Column"
$ws.Range("D2").Value = 5
$ws.Range("D3").Value = 5
$ws.Range("D4").Value = 5

$ws.Range("D2:D4").Select()
